$wb = $excel.ActiveWorkbook

# --- Sheet1: "Attendance Roster" ---
# Mark attendance ("P") for learners in column C, rows 12-26, skipping
# the rows that stay blank (16, 17, 20, 24).
$ws1 = $wb.Worksheets.Item("Sheet1")

$presentRows = @(12, 13, 14, 15, 18, 19, 21, 22, 23, 25, 26)
foreach ($r in $presentRows) {
    $ws1.Range("C$r").Value = "P"
}

# --- Sheet2: leave data intact, only cursor/tab state changes below ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# User finishes work on Sheet2 (cursor ends at F27) then switches back to
# Sheet1 (cursor ends at F13), which becomes the active tab on save.
$ws2.Range("F27").Select() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("F13").Select() | Out-Null
